$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7577165961265564
$ws.Range("B1").Value = 0.7929579615592957
$ws.Range("C1").Value = 2.777397394180298
$ws.Range("D1").Value = 5.191154479980469
$ws.Range("E1").Value = 1.188583970069885
